$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B1: "Cities" -> "Hobby"
$ws.Range("B1").Value = "Hobby"

# A2: "99004361" -> "99004359" (kept as text, matching the original inline-string type)
$ws.Range("A2").Value = "'99004359"

# B2: multi-line city list -> "Gardening"
$ws.Range("B2").Value = "Gardening"
